# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Sun Aug 11 20:35:25 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.612.35"
$ws.Range("E2").Value = "  -3.66%  "

$ws.Range("D3").Value = "2.560.54"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.03"
$ws.Range("E5").Value = "  -3.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.23"
$ws.Range("E6").Value = "  -5.70%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").Value = "  -6.20%  "

$ws.Range("D9").Value = "2.561.32"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.20"
$ws.Range("E10").Value = "  -7.59%  "

$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("E12").Value = "  -4.32%  "

$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").Value = "3.012.98"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "58.632.14"
$ws.Range("E15").Value = "  -3.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.55"
$ws.Range("E16").Value = "  -4.93%  "

$ws.Range("E17").Value = "  -4.80%  "

$ws.Range("D18").Value = "2.569.18"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  -4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.92"
$ws.Range("E20").Value = "  -5.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.09"
$ws.Range("E21").Value = "  -4.50%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -4.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.64"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -4.45%  "

$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("E27").Value = "  -6.68%  "

$ws.Range("D28").Value = "0.0₃0776"
$ws.Range("E28").Value = "  -8.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.84"
$ws.Range("E29").Value = "  -6.99%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").Value = "  -7.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.57"
$ws.Range("E32").Value = "  -3.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.01"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -6.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("E37").Value = "  -7.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.84"
$ws.Range("E38").Value = "  -1.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.53"
$ws.Range("E40").Value = "  -7.01%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -8.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "280.29"
$ws.Range("E42").Value = "  -2.70%  "

$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0983"
$ws.Range("E44").Value = "  -3.03%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.60"
$ws.Range("E48").Value = "  -4.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"

$ws.Range("D50").Value = "1.911.77"
$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("E51").Value = "  -7.60%  "
